$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RPM")

# Move the existing "K=" label and its helper formula one column to the right
# (K6 -> L6, L6 -> M6) to make room for the new I/J/K ratio columns.
$ws.Range("M6").Formula = "=(D10-D9)/50"
$ws.Range("M6").Style = "Normal"
$ws.Range("L6").Value = "K="
$ws.Range("K6").ClearContents() | Out-Null

# New ratio columns: I = D/A, J = G/A, K = AVERAGE(I:J)
$ws.Range("I5").Formula = "=D5/A5"
$ws.Range("J5").Formula = "=G5/A5"
$ws.Range("K5").Formula = "=AVERAGE(I5:J5)"

$ws.Range("I6:I11").Formula = "=D6/A6"
$ws.Range("J6:J11").Formula = "=G6/A6"
$ws.Range("K6:K12").Formula = "=AVERAGE(I6:J6)"

$ws.Range("I12").Formula = "=D12/A12"
$ws.Range("J12").Formula = "=G12/A12"

$ws.Range("I5:K12").Style = "Normal"

$ws.Range("K13").Formula = "=AVERAGE(K5:K12)"

# View changes: zoom + selection
$ws.Activate()
$excel.ActiveWindow.Zoom = 130
$ws.Range("K13").Select()

$wb.Save()
